$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 17 carries the footnote text that used to live in row 18 ("(1: ...)").
# Copy formatting (style s="9", matching rows 18/19) from the existing A18 cell
# so no new style/font entries get minted in styles.xml.
$ws.Range("A18").Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$ws.Range("A17").Value = "(1: Goo Goo Dolls, 1998.)"

# Row 18 now shows the footnote text that used to live in row 19 ("(2: ...)").
$ws.Range("A18").Value = "(2: I've got 99 problems and taxonomy is one.)"

# Row 19 gets a brand new annotation string.
$ws.Range("A19").Value = "(This note hath no reference.)"

# Restore the user's active selection to D7 (as captured in the saved sheetView).
$ws.Range("D7").Select() | Out-Null
